$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Merge the two runs "Generate BST" + " from postorder" into a
#    single run "Generate BST from postorder" (Find/Replace causes
#    the engine to coalesce the run).
# ------------------------------------------------------------------
$d.Content.Find.Execute("Generate BST from postorder", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Generate BST from postorder", 2)

# ------------------------------------------------------------------
# 2. Insert two new bulleted list paragraphs after it:
#       "Adjacency List"
#       "Compact List  "
#    New paragraphs created via InsertParagraphAfter() naturally
#    inherit the ListParagraph / numPr formatting of paragraph 10.
# ------------------------------------------------------------------
$pGenerateBst = $d.Paragraphs.Item(10)
$rEnd = $pGenerateBst.Range
$rEnd.Collapse(0)
$rEnd.InsertParagraphAfter()

$pAdjacency = $d.Paragraphs.Item(11)
$pAdjacency.Range.InsertAfter("Adjacency List")

$rAdjEnd = $pAdjacency.Range
$rAdjEnd.Collapse(0)
$rAdjEnd.InsertParagraphAfter()

$pCompact = $d.Paragraphs.Item(12)
$pCompact.Range.InsertAfter("Compact List  ")
# temporary placeholder character so the bookmark can be anchored
# without landing exactly on the paragraph-end boundary
$pCompact.Range.InsertAfter("X")

# ------------------------------------------------------------------
# 3. Move the "_GoBack" bookmark so it sits at the very end of the
#    new "Compact List  " paragraph (immediately before the
#    paragraph mark), matching the target layout.
# ------------------------------------------------------------------
$pCompactNow = $d.Paragraphs.Item(12)
$placeholderStart = $pCompactNow.Range.End - 2   # position right before the "X"

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$bmRange = $d.Range($placeholderStart, $placeholderStart)
$d.Bookmarks.Add("_GoBack", $bmRange)

# remove the temporary placeholder character
$xRange = $d.Range($placeholderStart, $placeholderStart + 1)
$xRange.Delete()

# ------------------------------------------------------------------
# 4. Remove two of the three trailing empty ListParagraph paragraphs,
#    leaving just one before the section properties.
# ------------------------------------------------------------------
$d.Paragraphs.Item(13).Range.Delete()
$d.Paragraphs.Item(13).Range.Delete()
